# Finish the "thinking question" in the pipeline-design workbook:
# extend the forwarding/hazard table (rows 39-45) with the MEMtoEX /
# WBtoEX / MEMtoID forwarding-path labels, and add a new data row for a
# second "store" case that also needs rt-forwarding.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at 45 (pushes the old row45.. down by one) ---
$ws.Rows.Item(45).Insert()

# The M39:P44 merged banner needs to grow to M39:P45 to cover the new row.
$ws.Range("M39:P45").Merge()
$ws.Range("M39:P45").HorizontalAlignment = -4108
$ws.Range("M39:P45").VerticalAlignment = -4108

# --- 2. Fill in the newly-revealed forwarding-path cells (rows 39-44) ---

# Row 39 (beq): MEM/WB -> ID/EX forwarding labelled MEMtoID
$ws.Cells.Item(39, 9).Value = "MEMtoID"   # I39
$ws.Cells.Item(39, 10).Value = "MEMtoID"  # J39
$ws.Cells.Item(39, 12).Value = "MEMtoID"  # L39

# Rows 40-43 (cal_r / cal_i / load / store): EX/MEM->EX (MEMtoEX) and
# MEM/WB->EX (WBtoEX) forwarding paths
40..43 | ForEach-Object {
    $r = $_
    $ws.Cells.Item($r, 5).Value  = "MEMtoEX"  # E
    $ws.Cells.Item($r, 6).Value  = "MEMtoEX"  # F
    $ws.Cells.Item($r, 8).Value  = "MEMtoEX"  # H
    $ws.Cells.Item($r, 9).Value  = "WBtoEX"   # I
    $ws.Cells.Item($r, 10).Value = "WBtoEX"   # J
    $ws.Cells.Item($r, 11).Value = "WBtoEX"   # K
    $ws.Cells.Item($r, 12).Value = "WBtoEX"   # L
}

# Row 44 (jr): same MEM/WB -> ID/EX forwarding as row 39
$ws.Cells.Item(44, 9).Value  = "MEMtoID"  # I44
$ws.Cells.Item(44, 10).Value = "MEMtoID"  # J44
$ws.Cells.Item(44, 12).Value = "MEMtoID"  # L44

# --- 3. Populate the new row 45: a second "store" case needing rt forwarding ---
$ws.Cells.Item(45, 2).Value  = "store"    # B45
$ws.Cells.Item(45, 3).Value  = "rt"       # C45
$ws.Cells.Item(45, 4).Value  = 2          # D45
$ws.Cells.Item(45, 5).Value  = "MEMtoEX"  # E45
$ws.Cells.Item(45, 6).Value  = "MEMtoEX"  # F45
$ws.Cells.Item(45, 7).Value  = "WBtoMEM"  # G45
$ws.Cells.Item(45, 8).Value  = "MEMtoEX"  # H45
$ws.Cells.Item(45, 9).Value  = "MEMtoID"  # I45
$ws.Cells.Item(45, 10).Value = "MEMtoID"  # J45
$ws.Cells.Item(45, 11).Value = "WBtoEX"   # K45
$ws.Cells.Item(45, 12).Value = "MEMtoID"  # L45

# Match the alignment used by the rest of the table (B:L, center-aligned)
$ws.Range("B45:L45").HorizontalAlignment = -4108

Write-Output "ok"
